# Adds 6 new species-observation rows (rows 3-8) to the Artfynd sheet,
# matching the source data for project "Kustpaketet" (reported 2025-06-24).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a literal text string (even if it
# looks like a number or a date, e.g. "1", "2025-06-24", "14:39") without
# leaving a non-default number format behind on the cell.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Helper: write a plain numeric value.
function Set-NumberCell($range, $value) {
    $range.Value = $value
}

# Helper: write a boolean (TRUE/FALSE) value.
function Set-BoolCell($range, $value) {
    $range.Value = $value
}

# ---- Row 3 ----
Set-NumberCell $ws.Range('A3') 131106808
Set-NumberCell $ws.Range('B3') 79000
Set-TextCell $ws.Range('D3') 'NT'
Set-NumberCell $ws.Range('E3') 6446
Set-TextCell $ws.Range('F3') 'Kolflarnlav'
Set-TextCell $ws.Range('G3') 'Carbonicola anthracophila'
Set-TextCell $ws.Range('H3') '(Nyl.) Bendiksby & Timdal'
Set-TextCell $ws.Range('I3') '1'
Set-TextCell $ws.Range('J3') 'dm²'
Set-TextCell $ws.Range('P3') 'Paljack, Mpd'
Set-NumberCell $ws.Range('Q3') 600356
Set-NumberCell $ws.Range('R3') 6973038
Set-NumberCell $ws.Range('S3') 10
Set-TextCell $ws.Range('T3') 'Västernorrland'
Set-TextCell $ws.Range('U3') 'Sundsvall'
Set-TextCell $ws.Range('V3') 'Medelpad'
Set-TextCell $ws.Range('W3') 'Liden'
Set-TextCell $ws.Range('X3') '2025_0368'
Set-TextCell $ws.Range('Y3') '2025-06-24'
Set-TextCell $ws.Range('Z3') '14:39'
Set-TextCell $ws.Range('AA3') '2025-06-24'
Set-TextCell $ws.Range('AB3') '14:39'
Set-BoolCell $ws.Range('AD3') $FALSE
Set-BoolCell $ws.Range('AE3') $FALSE
Set-BoolCell $ws.Range('AG3') $FALSE
Set-TextCell $ws.Range('AW3') 'David Isaksson'
Set-TextCell $ws.Range('AX3') 'David Isaksson'
Set-TextCell $ws.Range('AY3') 'Kustpaketet'

# ---- Row 4 ----
Set-NumberCell $ws.Range('A4') 131106807
Set-NumberCell $ws.Range('B4') 79862
Set-TextCell $ws.Range('D4') 'NT'
Set-NumberCell $ws.Range('E4') 6453
Set-TextCell $ws.Range('F4') 'Vedskivlav'
Set-TextCell $ws.Range('G4') 'Hertelidea botryosa'
Set-TextCell $ws.Range('H4') '(Fr.) Printzen & Kantvilas'
Set-TextCell $ws.Range('I4') '2'
Set-TextCell $ws.Range('J4') 'dm²'
Set-TextCell $ws.Range('P4') 'Paljack, Mpd'
Set-NumberCell $ws.Range('Q4') 600345
Set-NumberCell $ws.Range('R4') 6973042
Set-NumberCell $ws.Range('S4') 10
Set-TextCell $ws.Range('T4') 'Västernorrland'
Set-TextCell $ws.Range('U4') 'Sundsvall'
Set-TextCell $ws.Range('V4') 'Medelpad'
Set-TextCell $ws.Range('W4') 'Liden'
Set-TextCell $ws.Range('X4') '2025_0369'
Set-TextCell $ws.Range('Y4') '2025-06-24'
Set-TextCell $ws.Range('Z4') '14:40'
Set-TextCell $ws.Range('AA4') '2025-06-24'
Set-TextCell $ws.Range('AB4') '14:40'
Set-BoolCell $ws.Range('AD4') $FALSE
Set-BoolCell $ws.Range('AE4') $FALSE
Set-BoolCell $ws.Range('AG4') $FALSE
Set-TextCell $ws.Range('AW4') 'David Isaksson'
Set-TextCell $ws.Range('AX4') 'David Isaksson'
Set-TextCell $ws.Range('AY4') 'Kustpaketet'

# ---- Row 5 ----
Set-NumberCell $ws.Range('A5') 131106804
Set-NumberCell $ws.Range('B5') 79243
Set-TextCell $ws.Range('D5') 'NT'
Set-NumberCell $ws.Range('E5') 6425
Set-TextCell $ws.Range('F5') 'Garnlav'
Set-TextCell $ws.Range('G5') 'Alectoria sarmentosa'
Set-TextCell $ws.Range('H5') '(Ach.) Ach.'
Set-TextCell $ws.Range('P5') 'Paljack, Mpd'
Set-NumberCell $ws.Range('Q5') 600308
Set-NumberCell $ws.Range('R5') 6972996
Set-NumberCell $ws.Range('S5') 10
Set-TextCell $ws.Range('T5') 'Västernorrland'
Set-TextCell $ws.Range('U5') 'Sundsvall'
Set-TextCell $ws.Range('V5') 'Medelpad'
Set-TextCell $ws.Range('W5') 'Liden'
Set-TextCell $ws.Range('X5') '2025_0372'
Set-TextCell $ws.Range('Y5') '2025-06-24'
Set-TextCell $ws.Range('Z5') '14:53'
Set-TextCell $ws.Range('AA5') '2025-06-24'
Set-TextCell $ws.Range('AB5') '14:53'
Set-BoolCell $ws.Range('AD5') $FALSE
Set-BoolCell $ws.Range('AE5') $FALSE
Set-BoolCell $ws.Range('AG5') $FALSE
Set-TextCell $ws.Range('AW5') 'David Isaksson'
Set-TextCell $ws.Range('AX5') 'David Isaksson'
Set-TextCell $ws.Range('AY5') 'Kustpaketet'

# ---- Row 6 ----
Set-NumberCell $ws.Range('A6') 131106805
Set-NumberCell $ws.Range('B6') 79000
Set-TextCell $ws.Range('D6') 'NT'
Set-NumberCell $ws.Range('E6') 6446
Set-TextCell $ws.Range('F6') 'Kolflarnlav'
Set-TextCell $ws.Range('G6') 'Carbonicola anthracophila'
Set-TextCell $ws.Range('H6') '(Nyl.) Bendiksby & Timdal'
Set-TextCell $ws.Range('I6') '1'
Set-TextCell $ws.Range('J6') 'dm²'
Set-TextCell $ws.Range('P6') 'Paljack, Mpd'
Set-NumberCell $ws.Range('Q6') 600334
Set-NumberCell $ws.Range('R6') 6973001
Set-NumberCell $ws.Range('S6') 10
Set-TextCell $ws.Range('T6') 'Västernorrland'
Set-TextCell $ws.Range('U6') 'Sundsvall'
Set-TextCell $ws.Range('V6') 'Medelpad'
Set-TextCell $ws.Range('W6') 'Liden'
Set-TextCell $ws.Range('X6') '2025_0371'
Set-TextCell $ws.Range('Y6') '2025-06-24'
Set-TextCell $ws.Range('Z6') '14:50'
Set-TextCell $ws.Range('AA6') '2025-06-24'
Set-TextCell $ws.Range('AB6') '14:50'
Set-BoolCell $ws.Range('AD6') $FALSE
Set-BoolCell $ws.Range('AE6') $FALSE
Set-BoolCell $ws.Range('AG6') $FALSE
Set-TextCell $ws.Range('AW6') 'David Isaksson'
Set-TextCell $ws.Range('AX6') 'David Isaksson'
Set-TextCell $ws.Range('AY6') 'Kustpaketet'

# ---- Row 7 ----
Set-NumberCell $ws.Range('A7') 131106806
Set-NumberCell $ws.Range('B7') 79833
Set-TextCell $ws.Range('D7') 'NT'
Set-NumberCell $ws.Range('E7') 229821
Set-TextCell $ws.Range('F7') 'Vedflamlav'
Set-TextCell $ws.Range('G7') 'Ramboldia elabens'
Set-TextCell $ws.Range('H7') '(Fr.) Kantvilas & Elix'
Set-TextCell $ws.Range('I7') '2'
Set-TextCell $ws.Range('J7') 'cm²'
Set-TextCell $ws.Range('P7') 'Paljack, Mpd'
Set-NumberCell $ws.Range('Q7') 600353
Set-NumberCell $ws.Range('R7') 6973046
Set-NumberCell $ws.Range('S7') 10
Set-TextCell $ws.Range('T7') 'Västernorrland'
Set-TextCell $ws.Range('U7') 'Sundsvall'
Set-TextCell $ws.Range('V7') 'Medelpad'
Set-TextCell $ws.Range('W7') 'Liden'
Set-TextCell $ws.Range('X7') '2025_0370'
Set-TextCell $ws.Range('Y7') '2025-06-24'
Set-TextCell $ws.Range('Z7') '14:43'
Set-TextCell $ws.Range('AA7') '2025-06-24'
Set-TextCell $ws.Range('AB7') '14:43'
Set-BoolCell $ws.Range('AD7') $FALSE
Set-BoolCell $ws.Range('AE7') $FALSE
Set-BoolCell $ws.Range('AG7') $FALSE
Set-TextCell $ws.Range('AW7') 'David Isaksson'
Set-TextCell $ws.Range('AX7') 'David Isaksson'
Set-TextCell $ws.Range('AY7') 'Kustpaketet'

# ---- Row 8 ----
Set-NumberCell $ws.Range('A8') 131106803
Set-NumberCell $ws.Range('B8') 79000
Set-TextCell $ws.Range('D8') 'NT'
Set-NumberCell $ws.Range('E8') 6446
Set-TextCell $ws.Range('F8') 'Kolflarnlav'
Set-TextCell $ws.Range('G8') 'Carbonicola anthracophila'
Set-TextCell $ws.Range('H8') '(Nyl.) Bendiksby & Timdal'
Set-TextCell $ws.Range('I8') '1'
Set-TextCell $ws.Range('J8') 'm²'
Set-TextCell $ws.Range('P8') 'Paljack, Mpd'
Set-NumberCell $ws.Range('Q8') 600279
Set-NumberCell $ws.Range('R8') 6972994
Set-NumberCell $ws.Range('S8') 10
Set-TextCell $ws.Range('T8') 'Västernorrland'
Set-TextCell $ws.Range('U8') 'Sundsvall'
Set-TextCell $ws.Range('V8') 'Medelpad'
Set-TextCell $ws.Range('W8') 'Liden'
Set-TextCell $ws.Range('X8') '2025_0373'
Set-TextCell $ws.Range('Y8') '2025-06-24'
Set-TextCell $ws.Range('Z8') '14:59'
Set-TextCell $ws.Range('AA8') '2025-06-24'
Set-TextCell $ws.Range('AB8') '14:59'
Set-BoolCell $ws.Range('AD8') $FALSE
Set-BoolCell $ws.Range('AE8') $FALSE
Set-BoolCell $ws.Range('AG8') $FALSE
Set-TextCell $ws.Range('AW8') 'David Isaksson'
Set-TextCell $ws.Range('AX8') 'David Isaksson'
Set-TextCell $ws.Range('AY8') 'Kustpaketet'

